# Committing to new branch.
#
# Adds three new worksheets to the workbook:
#   - CypherOutput_Message  (copy of the existing "Message" sheet)
#   - StatOutput            (small 4-column stat summary table)
#   - StatOutput_Message    (two stacked "Message" blocks, second one
#                            carrying the new StatOutput cypher query)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Reusable "Message" style block values (Neo4j_URL / User_name / PWD / ...)
# ---------------------------------------------------------------------
$neo4jUrlLabel = 'Neo4j_URL:'
$neo4jUrlValue = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$userLabel = 'User_name:'
$userValue = 'neo4j'
$pwdLabel = 'PWD:'
$pwdValue = 'icdcDBneo4j0'
$cypherLabel = 'Cypher:'
$outputLabel = 'Output:'
$outputValue = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC39_Canine_Filter_Breed-WestHlnd_Neo4jData.xlsx'

# The original CypherOutput query (already present as a shared string).
$cypherOutputQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN [''West Highland White Terrier''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'

# The new StatOutput query used on the StatOutput_Message sheet.
$statOutputQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''West Highland White Terrier'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

function Write-MessageBlock($ws, $startRow, $cypherValue) {
    $ws.Cells.Item($startRow, 1).Value = $neo4jUrlLabel
    $ws.Cells.Item($startRow + 1, 1).Value = $neo4jUrlValue
    $ws.Cells.Item($startRow + 2, 1).Value = $userLabel
    $ws.Cells.Item($startRow + 3, 1).Value = $userValue
    $ws.Cells.Item($startRow + 4, 1).Value = $pwdLabel
    $ws.Cells.Item($startRow + 5, 1).Value = $pwdValue
    $ws.Cells.Item($startRow + 6, 1).Value = $cypherLabel
    $ws.Cells.Item($startRow + 7, 1).Value = $cypherValue
    $ws.Cells.Item($startRow + 8, 1).Value = $outputLabel
    $ws.Cells.Item($startRow + 9, 1).Value = $outputValue
}

# ---------------------------------------------------------------------
# Sheet: CypherOutput_Message (copy of Message)
# ---------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"
Write-MessageBlock $cypherOutputMessage 1 $cypherOutputQuery

# ---------------------------------------------------------------------
# Sheet: StatOutput
# ---------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"
$statOutput.Cells.Item(1, 1).Value = "number_of_files"
$statOutput.Cells.Item(1, 2).Value = "number_of_sample"
$statOutput.Cells.Item(1, 3).Value = "number_of_cases"
$statOutput.Cells.Item(1, 4).Value = "number_of_study"
$statOutput.Cells.Item(2, 1).Value = "'1"
$statOutput.Cells.Item(2, 2).Value = "'2"
$statOutput.Cells.Item(2, 3).Value = "'1"
$statOutput.Cells.Item(2, 4).Value = "'1"

# ---------------------------------------------------------------------
# Sheet: StatOutput_Message (two stacked Message blocks; second block
# carries the new StatOutput cypher text)
# ---------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$statOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$statOutputMessage.Name = "StatOutput_Message"
Write-MessageBlock $statOutputMessage 1 $cypherOutputQuery
Write-MessageBlock $statOutputMessage 11 $statOutputQuery

# ---------------------------------------------------------------------
# Restore CypherOutput as the active/selected tab (matches original state).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
